$wb = $excel.ActiveWorkbook

# --- Sheet "对公业务台账202312": update column T (三个月发生数) ---
$wsT = $wb.Worksheets.Item("对公业务台账202312")
$cell = $wsT.Range("T2")
$cell.Value = "'46.00"
$cell.Style = "Normal"
$cell = $wsT.Range("T3")
$cell.Value = "'11.00"
$cell.Style = "Normal"
$cell = $wsT.Range("T4")
$cell.Value = "'2.00"
$cell.Style = "Normal"
$cell = $wsT.Range("T5")
$cell.Value = "'25.00"
$cell.Style = "Normal"
$cell = $wsT.Range("T7")
$cell.Value = "'161.00"
$cell.Style = "Normal"
$cell = $wsT.Range("T9")
$cell.Value = "'45.00"
$cell.Style = "Normal"
$cell = $wsT.Range("T10")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$cell = $wsT.Range("T11")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$cell = $wsT.Range("T12")
$cell.Value = "'292.00"
$cell.Style = "Normal"

# --- Sheet "供应链放还款202312": update column H (今年日均) ---
$wsH = $wb.Worksheets.Item("供应链放还款202312")
$wsH.Range("H2").Value = 5317.36
$wsH.Range("H3").Value = 1108.05
$wsH.Range("H4").Value = 6530.42
$wsH.Range("H5").Value = 191.05
$wsH.Range("H6").Value = 85.85
$wsH.Range("H7").Value = 20585.53
$wsH.Range("H8").Value = 82773.94
$wsH.Range("H9").Value = 100.04
$wsH.Range("H10").Value = 7.78
$wsH.Range("H11").Value = 8.98
$wsH.Range("H12").Value = 18169.4
$wsH.Range("H13").Value = 20526.27
$wsH.Range("H14").Value = 22230.52
$wsH.Range("H15").Value = 248625.26
$wsH.Range("H16").Value = 1588.19
$wsH.Range("H17").Value = 307.71
$wsH.Range("H18").Value = 428156.35
